$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the header formatting (bold, centered, bordered) from the
# existing "RF" header cell (E1) into the three new header cells F1:H1.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New header labels for the outlier-flag columns.
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Populate the new boolean outlier-flag columns (F, G, H) for every data
# row (2-18) with FALSE.
$ws.Range("F2:H18").Value = $false
